$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new value looks like a plain number as Text,
# so Excel keeps them as literal strings (matching the source data,
# which stores every Price value as text) instead of silently
# converting to a floating point number.
$textCells = @("D4", "D5", "D6", "D8", "D11", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D25", "D26", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "61.283.44"
$ws.Range("E2").Value = "  -4.00%  "

$ws.Range("D3").Value = "2.460.12"
$ws.Range("E3").Value = "  -6.59%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "547.22"
$ws.Range("E5").Value = "  -5.47%  "

$ws.Range("D6").Value = "146.38"
$ws.Range("E6").Value = "  -6.65%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -6.79%  "

$ws.Range("D9").Value = "2.458.53"
$ws.Range("E9").Value = "  -6.59%  "

$ws.Range("E10").Value = "  -10.13%  "

$ws.Range("D11").Value = "5.44"
$ws.Range("E11").Value = "  -6.65%  "

$ws.Range("E12").Value = "  -1.88%  "

$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -8.52%  "

$ws.Range("D14").Value = "26.08"
$ws.Range("E14").Value = "  -9.37%  "

$ws.Range("D15").Value = "2.898.94"
$ws.Range("E15").Value = "  -6.85%  "

$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -9.72%  "

$ws.Range("D17").Value = "61.178.34"
$ws.Range("E17").Value = "  -4.05%  "

$ws.Range("D18").Value = "2.456.67"
$ws.Range("E18").Value = "  -6.95%  "

$ws.Range("D19").Value = "11.11"
$ws.Range("E19").Value = "  -8.80%  "

$ws.Range("D20").Value = "7.05"
$ws.Range("E20").Value = "  -8.72%  "

$ws.Range("D21").Value = "4.17"
$ws.Range("E21").Value = "  -7.92%  "

$ws.Range("D22").Value = "317.95"
$ws.Range("E22").Value = "  -7.51%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").Value = "63.82"
$ws.Range("E25").Value = "  -6.67%  "

$ws.Range("B26").Value = "Bittensor"
$ws.Range("C26").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D26").Value = "550.76"
$ws.Range("E26").Value = "  -5.34%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.578.09"
$ws.Range("E27").Value = "  -6.91%  "

$ws.Range("D28").Value = "0.0₃0961"
$ws.Range("E28").Value = "  -14.65%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  -11.70%  "

$ws.Range("D31").Value = "8.22"
$ws.Range("E31").Value = "  -11.13%  "

$ws.Range("D32").Value = "7.60"
$ws.Range("E32").Value = "  -8.43%  "

$ws.Range("D33").Value = "0.147"
$ws.Range("E33").Value = "  -8.26%  "

$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  -8.41%  "

$ws.Range("E35").Value = "  -9.10%  "

$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  -12.29%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "4.81"
$ws.Range("E38").Value = "  -12.20%  "

$ws.Range("D39").Value = "0.378"
$ws.Range("E39").Value = "  -6.46%  "

$ws.Range("D40").Value = "18.41"
$ws.Range("E40").Value = "  -7.01%  "

$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  -7.33%  "

$ws.Range("D42").Value = "140.88"
$ws.Range("E42").Value = "  -8.66%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "40.41"
$ws.Range("E44").Value = "  -4.05%  "

$ws.Range("E45").Value = "  -11.02%  "

$ws.Range("D46").Value = "146.13"
$ws.Range("E46").Value = "  -10.55%  "

$ws.Range("D47").Value = "3.59"
$ws.Range("E47").Value = "  -8.56%  "

$ws.Range("D48").Value = "21.47"
$ws.Range("E48").Value = "  -11.12%  "

$ws.Range("D49").Value = "0.0535"
$ws.Range("E49").Value = "  -8.96%  "

$ws.Range("D50").Value = "0.587"
$ws.Range("E50").Value = "  -7.49%  "

$ws.Range("D51").Value = "0.0934"
$ws.Range("E51").Value = "  -6.99%  "

